$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '27.651.12'
$ws.Cells.Item(2, 5).Value = '  -0.47%  '

$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '1.899.42'
$ws.Cells.Item(3, 5).Value = '  +0.06%  '

$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.002'
$ws.Cells.Item(4, 5).Value = '  -0.24%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '311.47'
$ws.Cells.Item(5, 5).Value = '  -1.03%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '1.000'
$ws.Cells.Item(6, 5).Value = '  -0.41%  '

$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.5152'
$ws.Cells.Item(7, 5).Value = '  +6.93%  '

$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.3772'
$ws.Cells.Item(8, 5).Value = '  -0.93%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.07236'
$ws.Cells.Item(9, 5).Value = '  -1.40%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '20.90'
$ws.Cells.Item(10, 5).Value = '  +1.34%  '

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.8901'
$ws.Cells.Item(11, 5).Value = '  -3.70%  '

$ws.Cells.Item(12, 2).Value = 'TRON'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.07648'
$ws.Cells.Item(12, 5).Value = '  -0.90%  '

$ws.Cells.Item(13, 2).Value = 'WrappedEther'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '1.905.63'
$ws.Cells.Item(13, 5).Value = '  +0.15%  '

$ws.Cells.Item(14, 5).Value = '  -0.65%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '91.90'
$ws.Cells.Item(15, 5).Value = '  +0.50%  '

$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '1.003'
$ws.Cells.Item(16, 5).Value = '  -0.37%  '

$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '0.000008712'
$ws.Cells.Item(17, 5).Value = '  -1.22%  '

$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '1.001'
$ws.Cells.Item(18, 5).Value = '  -0.42%  '

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '27.702.24'
$ws.Cells.Item(19, 5).Value = '  -0.46%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '14.43'
$ws.Cells.Item(20, 5).Value = '  -0.86%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '5.141'
$ws.Cells.Item(21, 5).Value = '  +0.08%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '2.148.82'
$ws.Cells.Item(22, 5).Value = '  -1.36%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '10.81'
$ws.Cells.Item(23, 5).Value = '  -0.38%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '6.576'
$ws.Cells.Item(24, 5).Value = '  -0.42%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '153.58'
$ws.Cells.Item(25, 5).Value = '  -0.46%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '1.862'
$ws.Cells.Item(26, 5).Value = '  -2.84%  '

$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '2.180'
$ws.Cells.Item(27, 5).Value = '  +2.33%  '

$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '18.27'
$ws.Cells.Item(28, 5).Value = '  -0.80%  '

$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '114.67'
$ws.Cells.Item(29, 5).Value = '  -1.69%  '

$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '4.833'
$ws.Cells.Item(30, 5).Value = '  -2.11%  '

$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '0.08946'
$ws.Cells.Item(31, 5).Value = '  -0.15%  '

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '3.185'
$ws.Cells.Item(32, 5).Value = '  +0.34%  '

$ws.Cells.Item(33, 2).Value = 'Filecoin'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '4.787'
$ws.Cells.Item(33, 5).Value = '  +3.12%  '

$ws.Cells.Item(34, 2).Value = 'ARBITRUM'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '1.227'
$ws.Cells.Item(34, 5).Value = '  -1.37%  '

$ws.Cells.Item(35, 2).Value = 'ImmutableX'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.7761'
$ws.Cells.Item(35, 5).Value = '  +1.64%  '

$ws.Cells.Item(36, 2).Value = 'VeChain'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.02086'
$ws.Cells.Item(36, 5).Value = '  +2.28%  '

$ws.Cells.Item(37, 2).Value = 'RenderToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '2.605'
$ws.Cells.Item(37, 5).Value = '  +2.69%  '

$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '3.058'
$ws.Cells.Item(38, 5).Value = '  +2.33%  '

$ws.Cells.Item(39, 5).Value = '  -0.40%  '

$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.5478'
$ws.Cells.Item(40, 5).Value = '  +0.61%  '

$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.05258'
$ws.Cells.Item(41, 5).Value = '  -0.27%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '6.690'
$ws.Cells.Item(42, 5).Value = '  -3.74%  '

$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '113.12'
$ws.Cells.Item(43, 5).Value = '  +3.38%  '

$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '8.452'
$ws.Cells.Item(44, 5).Value = '  +1.61%  '

$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.1497'
$ws.Cells.Item(45, 5).Value = '  -1.39%  '

$ws.Cells.Item(46, 5).Value = '  -0.49%  '

$ws.Cells.Item(47, 5).Value = '  -1.78%  '

$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '0.9994'
$ws.Cells.Item(48, 5).Value = '  -0.51%  '

$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.609'
$ws.Cells.Item(49, 5).Value = '  -2.06%  '

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '66.55'
$ws.Cells.Item(50, 5).Value = '  -1.83%  '

$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.06002'
$ws.Cells.Item(51, 5).Value = '  -1.15%  '
